# Update "想去人数" (want-to-go count) figures in column F for both the
# "展览" sheet and the aggregated "全部类型" sheet, matching the refreshed
# scrape output referenced by commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F
$updates = @{
    3  = 8552
    4  = 1523
    6  = 393
    7  = 257
    9  = 32
    11 = 45
    13 = 1253
    14 = 240
    15 = 81
    16 = 143
    17 = 100
    18 = 130
    19 = 78
    20 = 118
    21 = 106
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
